# Updated files for crab catch estimation
# Add the 2023 winter CRC totals (row 18, columns D/E) that were previously
# missing, and update the active selection to reflect where the editor
# left off (H22), clearing the prior scrolled/frozen top-left cell state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# total.crc.winter for 2023
$ws.Range("D18").Value = 32975
# reported.crc.winter for 2023
$ws.Range("E18").Value = 15204

# Reflect the saved cursor/selection position from the edit session.
$ws.Activate()
$ws.Range("H22").Select()
